# The "Form Responses 1" sheet holds a column of (otherwise unlabeled)
# numeric offsets in column A (formatted as a date/time, hence the
# `#VALUE!` display for negative serials) and the respondents' first
# names in column B. The edit bumps every data row's offset in column A
# from -16 to -18, leaving the names in column B untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Form Responses 1")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 36 }

$ws.Range("A2:A$lastRow").Value = -18
